$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Qid values in column B for rows 2-6 (newly added iAuthor TC's)
$ws.Range("B2").Value = 622
$ws.Range("B3").Value = 639
$ws.Range("B4").Value = 640
$ws.Range("B5").Value = 641
$ws.Range("B6").Value = 642

# Row 4's cAnswer (column G) changes from "A" to "B"
$ws.Range("G4").Value = "B"
